$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) contain numeric-looking / percentage-looking
# text that must stay stored as literal text, so force a text number format
# before assigning the value (avoids Excel auto-converting "301.71" or "-2.60%"
# into numeric/percentage values).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '301.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.60%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.51%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.075'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.52%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07929'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-2.84%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.886'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-8.49%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.779'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.04%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9284'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.34%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1378'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '28.37%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1898'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.77%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09079'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.28%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03441'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.88%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09836'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.76%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001409'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.27%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005870'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.76%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.532'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.67%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.048'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.98%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.982'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.31%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3425'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.93%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.18%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.032'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.33%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.44%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.88%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001213'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.07%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004758'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.47%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-1.60%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-32.57%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01853'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-5.42%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04747'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.01%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007326'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.10%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009730'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.68%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.02%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002111'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.94%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01095'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.56%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006235'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-5.41%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.04%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.66'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-64.65%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '10.48%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.04%'
